$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before D (2018-12-31 and 2018-09-30 quarters),
# shifting the existing D:K data right to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: copy number formatting from F:G (the old D:E, now shifted) down onto
# the new D:E columns for all data rows, so the new columns pick up the same
# date / number styles as the rest of the table
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: populate the new D:E columns with the newly reported quarterly data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 15700
$ws.Range("E8").Value = 18800
$ws.Range("D9").Value = 9700
$ws.Range("E9").Value = 11300
$ws.Range("D10").Value = 6000
$ws.Range("E10").Value = 7500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 14200
$ws.Range("E17").Value = 16300
$ws.Range("D18").Value = 1500
$ws.Range("E18").Value = 2500
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 2300
$ws.Range("D21").Value = 1600
$ws.Range("E21").Value = 4900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1500
$ws.Range("E23").Value = 4900
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1600
$ws.Range("E26").Value = 4900
$ws.Range("D27").Value = 1300
$ws.Range("E27").Value = 4500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 1100
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = -2300
$ws.Range("D33").Value = 2500
$ws.Range("E33").Value = 4500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 2500
$ws.Range("E35").Value = 4500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 31300
$ws.Range("E41").Value = 34400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 25000
$ws.Range("E43").Value = 24100
$ws.Range("D44").Value = 3100
$ws.Range("E44").Value = 3900
$ws.Range("D45").Value = 3000
$ws.Range("E45").Value = 1600
$ws.Range("D46").Value = 62400
$ws.Range("E46").Value = 64000
$ws.Range("D47").Value = 2600
$ws.Range("E47").Value = 2500
$ws.Range("D48").Value = 89100
$ws.Range("E48").Value = 87800
$ws.Range("D49").Value = 9900
$ws.Range("E49").Value = 11100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8600
$ws.Range("E52").Value = 7200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 172500
$ws.Range("E54").Value = 172500
$ws.Range("D57").Value = 4600
$ws.Range("E57").Value = 6500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 1100
$ws.Range("D59").Value = 3300
$ws.Range("E59").Value = 1300
$ws.Range("D60").Value = 7900
$ws.Range("E60").Value = 8900
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 900
$ws.Range("E62").Value = 1500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 17500
$ws.Range("E66").Value = 19000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 59300
$ws.Range("E72").Value = 58100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 155000
$ws.Range("E76").Value = 153400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 2500
$ws.Range("E81").Value = 4500
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 3700
$ws.Range("E89").Value = 7500
$ws.Range("D91").Value = -500
$ws.Range("E91").Value = -2100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -2600
$ws.Range("E94").Value = -6700
$ws.Range("D96").Value = -1300
$ws.Range("E96").Value = -1300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2400
$ws.Range("E100").Value = -500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -1300
$ws.Range("E102").Value = 300

# Step 4: apply restatements to a handful of existing (shifted) cells whose
# historical values were revised in this update
$ws.Range("H8").Value = 11800
$ws.Range("H9").Value = 7400
$ws.Range("H10").Value = 4400
$ws.Range("H14").Value = 2900
$ws.Range("H17").Value = 12900
$ws.Range("H18").Value = -1100
$ws.Range("H20").Value = 700
$ws.Range("H21").Value = "NA"
$ws.Range("I21").Value = "NA"
$ws.Range("H23").Value = -400
$ws.Range("H26").Value = -500
$ws.Range("H27").Value = 100
$ws.Range("H29").Value = 1600
$ws.Range("H32").Value = -700
$ws.Range("H41").Value = 45500
$ws.Range("H43").Value = 16100
$ws.Range("H44").Value = 1800
$ws.Range("H45").Value = 3300
$ws.Range("H48").Value = 73000
$ws.Range("H49").Value = 11200
$ws.Range("H52").Value = 11000
$ws.Range("H57").Value = 3500
$ws.Range("H59").Value = 3400
$ws.Range("H89").Value = 1600
$ws.Range("I89").Value = 2900
$ws.Range("F91").Value = 2800
$ws.Range("H91").Value = -3900
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = -800
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
